$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 29.253501
$ws.Range("H2").Value = 87.760503
$ws.Range("I2").Value = 0.7876335333413836
$ws.Range("J2").Value = 0.7876335333413838
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.236283
$ws.Range("N2").Value = 0.708849
$ws.Range("O2").Value = 0.001461516295904947
$ws.Range("P2").Value = 0.001461516295904947
$ws.Range("Q2").Value = 6.912104976783
$ws.Range("R2").Value = 62.208944791047
$ws.Range("S2").Value = 0.001151139244179625
$ws.Range("T2").Value = 0.001151139244179625

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 29.253501
$ws.Range("H3").Value = 87.760503
$ws.Range("I3").Value = 0.7876335333413836
$ws.Range("J3").Value = 0.7876335333413838
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 160.8390196666666
$ws.Range("N3").Value = 482.517059
$ws.Range("O3").Value = 0.9948614511421033
$ws.Range("P3").Value = 0.9948614511421032
$ws.Range("Q3").Value = 4705.104422657852
$ws.Range("R3").Value = 42345.93980392067
$ws.Range("S3").Value = 0.7835862399481912
$ws.Range("T3").Value = 0.7835862399481912

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 29.253501
$ws.Range("H4").Value = 87.760503
$ws.Range("I4").Value = 0.7876335333413836
$ws.Range("J4").Value = 0.7876335333413838
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.119972
$ws.Range("N4").Value = 0.359916
$ws.Range("O4").Value = 0.0007420806111836584
$ws.Range("P4").Value = 0.0007420806111836582
$ws.Range("Q4").Value = 3.509601021972
$ws.Range("R4").Value = 31.586409197748
$ws.Range("S4").Value = 0.0005844875738107183
$ws.Range("T4").Value = 0.0005844875738107182

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 29.253501
$ws.Range("H5").Value = 87.760503
$ws.Range("I5").Value = 0.7876335333413836
$ws.Range("J5").Value = 0.7876335333413838
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.474493
$ws.Range("N5").Value = 1.423479
$ws.Range("O5").Value = 0.002934951950808252
$ws.Range("P5").Value = 0.002934951950808251
$ws.Range("Q5").Value = 13.880581449993
$ws.Range("R5").Value = 124.925233049937
$ws.Range("S5").Value = 0.00231166657520229
$ws.Range("T5").Value = 0.00231166657520229

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.749137666666667
$ws.Range("H6").Value = 14.247413
$ws.Range("I6").Value = 0.1278677748937237
$ws.Range("J6").Value = 0.1278677748937237
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.236283
$ws.Range("N6").Value = 0.708849
$ws.Range("O6").Value = 0.001461516295904947
$ws.Range("P6").Value = 0.001461516295904947
$ws.Range("Q6").Value = 1.122140495293
$ws.Range("R6").Value = 10.099264457637
$ws.Range("S6").Value = 0.0001868808367282827
$ws.Range("T6").Value = 0.0001868808367282827

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4.749137666666667
$ws.Range("H7").Value = 14.247413
$ws.Range("I7").Value = 0.1278677748937237
$ws.Range("J7").Value = 0.1278677748937237
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 160.8390196666666
$ws.Range("N7").Value = 482.517059
$ws.Range("O7").Value = 0.9948614511421033
$ws.Range("P7").Value = 0.9948614511421032
$ws.Range("Q7").Value = 763.8466465687073
$ws.Range("R7").Value = 6874.619819118368
$ws.Range("S7").Value = 0.1272107200850818
$ws.Range("T7").Value = 0.1272107200850818

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.749137666666667
$ws.Range("H8").Value = 14.247413
$ws.Range("I8").Value = 0.1278677748937237
$ws.Range("J8").Value = 0.1278677748937237
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.119972
$ws.Range("N8").Value = 0.359916
$ws.Range("O8").Value = 0.0007420806111836584
$ws.Range("P8").Value = 0.0007420806111836582
$ws.Range("Q8").Value = 0.5697635441453334
$ws.Range("R8").Value = 5.127871897308001
$ws.Range("S8").Value = 0.00009488819654382893
$ws.Range("T8").Value = 0.00009488819654382892

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.749137666666667
$ws.Range("H9").Value = 14.247413
$ws.Range("I9").Value = 0.1278677748937237
$ws.Range("J9").Value = 0.1278677748937237
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.474493
$ws.Range("N9").Value = 1.423479
$ws.Range("O9").Value = 0.002934951950808252
$ws.Range("P9").Value = 0.002934951950808251
$ws.Range("Q9").Value = 2.253432578869667
$ws.Range("R9").Value = 20.280893209827
$ws.Range("S9").Value = 0.0003752857753698448
$ws.Range("T9").Value = 0.0003752857753698448

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.138366333333334
$ws.Range("H10").Value = 9.415099000000001
$ws.Range("I10").Value = 0.08449869176489255
$ws.Range("J10").Value = 0.08449869176489258
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.236283
$ws.Range("N10").Value = 0.708849
$ws.Range("O10").Value = 0.001461516295904947
$ws.Range("P10").Value = 0.001461516295904947
$ws.Range("Q10").Value = 0.7415426123390001
$ws.Range("R10").Value = 6.673883511051001
$ws.Range("S10").Value = 0.0001234962149970396
$ws.Range("T10").Value = 0.0001234962149970397

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.138366333333334
$ws.Range("H11").Value = 9.415099000000001
$ws.Range("I11").Value = 0.08449869176489255
$ws.Range("J11").Value = 0.08449869176489258
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 160.8390196666666
$ws.Range("N11").Value = 482.517059
$ws.Range("O11").Value = 0.9948614511421033
$ws.Range("P11").Value = 0.9948614511421032
$ws.Range("Q11").Value = 504.7717644082045
$ws.Range("R11").Value = 4542.945879673842
$ws.Range("S11").Value = 0.0840644911088303
$ws.Range("T11").Value = 0.08406449110883032

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.138366333333334
$ws.Range("H12").Value = 9.415099000000001
$ws.Range("I12").Value = 0.08449869176489255
$ws.Range("J12").Value = 0.08449869176489258
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.119972
$ws.Range("N12").Value = 0.359916
$ws.Range("O12").Value = 0.0007420806111836584
$ws.Range("P12").Value = 0.0007420806111836582
$ws.Range("Q12").Value = 0.3765160857426668
$ws.Range("R12").Value = 3.388644771684
$ws.Range("S12").Value = 0.00006270484082911102
$ws.Range("T12").Value = 0.00006270484082911102

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.138366333333334
$ws.Range("H13").Value = 9.415099000000001
$ws.Range("I13").Value = 0.08449869176489255
$ws.Range("J13").Value = 0.08449869176489258
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.474493
$ws.Range("N13").Value = 1.423479
$ws.Range("O13").Value = 0.002934951950808252
$ws.Range("P13").Value = 0.002934951950808251
$ws.Range("Q13").Value = 1.489132856602334
$ws.Range("R13").Value = 13.402195709421
$ws.Range("S13").Value = 0.0002479996002361166
$ws.Range("T13").Value = 0.0002479996002361166
